# Daily attendance processing - 2025-10-27 20:20:58
# Reorders the "Recorded By" names/emails in column G for specific rows
# on the "Session Analysis Results" sheet (values are re-sequenced, same set).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "G2"; Value = "backup@backdoor.com, system, System" },
    @{ Cell = "G3"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G6"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G7"; Value = "System, admin@admin.com" },
    @{ Cell = "G10"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G11"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G12"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G13"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G14"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G15"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G17"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G18"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G19"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G20"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G21"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G22"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G24"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G29"; Value = "backup@backdoor.com, system, System" },
    @{ Cell = "G30"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G33"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G34"; Value = "System, admin@admin.com" },
    @{ Cell = "G37"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G38"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G39"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G40"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G41"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G42"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G44"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G45"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G46"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G47"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G48"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G49"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G51"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G56"; Value = "backup@backdoor.com, system, System" },
    @{ Cell = "G57"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G60"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G61"; Value = "System, admin@admin.com" },
    @{ Cell = "G64"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G65"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G66"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G67"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G68"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G69"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G71"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G72"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G73"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G74"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G75"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G76"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G78"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G86"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G87"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G88"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G89"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G90"; Value = "admin@admin.com, dnasr281@gmail.com" },
    @{ Cell = "G93"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G95"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G96"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G97"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G99"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G102"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G112"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G113"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G114"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G115"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G116"; Value = "admin@admin.com, dnasr281@gmail.com" },
    @{ Cell = "G119"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G121"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G122"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G123"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G125"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G128"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G138"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G139"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G140"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G141"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G142"; Value = "admin@admin.com, dnasr281@gmail.com" },
    @{ Cell = "G145"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G147"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G148"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G149"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G151"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G154"; Value = "System, dnasr281@gmail.com" }
)

foreach ($change in $changes) {
    $ws.Range($change.Cell).Value = $change.Value
}
